# Adds a "style" property row to the existing food-group sheets, and adds a
# brand new "noodles" sheet (inserted right before "bread") describing a
# csv-imported food item.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a cell value as a literal TEXT string even when the text
# looks like a boolean keyword ("true"/"false") which Excel would otherwise
# auto-convert to a native Boolean on entry. We build it as a formula that
# evaluates to the desired string, then convert the formula to a static
# value via copy / paste-special (values only) -- this keeps the cell typed
# as a shared-string without leaving any quote-prefix / text-format style
# behind on the cell.
# ---------------------------------------------------------------------------
function Set-TextValue {
    param($range, [string]$text)
    $escaped = $text.Replace("""", """""")
    $range.Formula = "=""" + $escaped + """"
    $range.Copy() | Out-Null
    $range.PasteSpecial(-4163) | Out-Null   # xlPasteValues
}

# ---- carbs: add "style" = "pass" --------------------------------------
$ws = $wb.Worksheets.Item("carbs")
$ws.Range("A8").Value = "style"
$ws.Range("B8").Value = "pass"
$ws.Range("B8").Select() | Out-Null

# ---- rice: add "style" = "asian" ---------------------------------------
$ws = $wb.Worksheets.Item("rice")
$ws.Range("A6").Value = "style"
$ws.Range("B6").Value = "asian"
$ws.Range("C31").Select() | Out-Null

# ---- potatoes: add "style" = "western" ---------------------------------
$ws = $wb.Worksheets.Item("potatoes")
$ws.Range("A6").Value = "style"
$ws.Range("B6").Value = "western"
$ws.Range("B6").Select() | Out-Null

# ---- oatmeal: add "style" = "western" -----------------------------------
$ws = $wb.Worksheets.Item("oatmeal")
$ws.Range("A6").Value = "style"
$ws.Range("B6").Value = "western"
$ws.Range("B6").Select() | Out-Null

# ---- grits: add "style" = "american" ------------------------------------
$ws = $wb.Worksheets.Item("grits")
$ws.Range("A6").Value = "style"
$ws.Range("B6").Value = "american"
$ws.Range("A7").Select() | Out-Null

# ---- grain: add "style" = "pass" -----------------------------------------
$ws = $wb.Worksheets.Item("grain")
$ws.Range("A5").Value = "style"
$ws.Range("B5").Value = "pass"
$ws.Range("B5").Select() | Out-Null

# ---- noodles: brand-new sheet, inserted right before "bread" -------------
$breadSheet = $wb.Worksheets.Item("bread")
$ws = $wb.Worksheets.Add($breadSheet)
$ws.Name = "noodles"

$ws.Range("A1").Value = "property"
$ws.Range("B1").Value = "value"

$ws.Range("A2").Value = "name"
$ws.Range("B2").Value = "noodles"

$ws.Range("A3").Value = "healthy"
Set-TextValue $ws.Range("B3") "true"

$ws.Range("A4").Value = "food super group"
$ws.Range("B4").Value = "carbs"

$ws.Range("A5").Value = "gluten-free substitute"
$ws.Range("B5").Value = $true

$ws.Range("A6").Value = "style"
$ws.Range("B6").Value = "asian"

$ws.Range("A7").Select() | Out-Null
